# Scheduled market-data sync for the Aegis_Profits workbook.
# Refreshes computed price/profit columns (H-N) per Leve row based on
# the latest Market Board averages pulled by the automation runner.
#   H = currentAveragePrice        K = LevePriceNQ
#   I = currentAveragePriceNQ      L = LevePriceHQ
#   J = currentAveragePriceHQ      M = LeveProfitNQ   N = LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 6500  # H13: 8000 -> 6500
$ws.Cells.Item(13, 10).Value = 6500  # J13: 8000 -> 6500
$ws.Cells.Item(13, 12).Value = 6500  # L13: 8000 -> 6500
$ws.Cells.Item(13, 14).Value = -6838  # N13: -8338 -> -6838
$ws.Cells.Item(33, 8).Value = 1047.8572  # H33: 1616.1666 -> 1047.8572
$ws.Cells.Item(33, 9).Value = 222.38889  # I33: 239.4 -> 222.38889
$ws.Cells.Item(33, 10).Value = 6000.6665  # J33: 8500 -> 6000.6665
$ws.Cells.Item(33, 11).Value = 222.38889  # K33: 239.4 -> 222.38889
$ws.Cells.Item(33, 12).Value = 6000.6665  # L33: 8500 -> 6000.6665
$ws.Cells.Item(33, 13).Value = 6.611109999999996  # M33: -10.40000000000001 -> 6.611109999999996
$ws.Cells.Item(33, 14).Value = -6458.6665  # N33: -8958 -> -6458.6665
$ws.Cells.Item(108, 8).Value = 32829.332  # H108: 39744 -> 32829.332
$ws.Cells.Item(108, 10).Value = 32829.332  # J108: 39744 -> 32829.332
$ws.Cells.Item(108, 12).Value = 32829.332  # L108: 39744 -> 32829.332
$ws.Cells.Item(108, 14).Value = -40509.332  # N108: -47424 -> -40509.332
$ws.Cells.Item(112, 8).Value = 1146.4445  # H112: 1167.091 -> 1146.4445
$ws.Cells.Item(112, 10).Value = 1339.7142  # J112: 1322 -> 1339.7142
$ws.Cells.Item(112, 12).Value = 4019.1426  # L112: 3966 -> 4019.1426
$ws.Cells.Item(112, 14).Value = -6235.142599999999  # N112: -6182 -> -6235.142599999999
$ws.Cells.Item(129, 8).Value = 3354.513  # H129: 3240.2195 -> 3354.513
$ws.Cells.Item(129, 9).Value = 5718.421  # I129: 5718.9473 -> 5718.421
$ws.Cells.Item(129, 10).Value = 1108.8  # J129: 1099.5 -> 1108.8
$ws.Cells.Item(129, 11).Value = 17155.263  # K129: 17156.8419 -> 17155.263
$ws.Cells.Item(129, 12).Value = 3326.4  # L129: 3298.5 -> 3326.4
$ws.Cells.Item(129, 13).Value = -12155.263  # M129: -12156.8419 -> -12155.263
$ws.Cells.Item(129, 14).Value = -13326.4  # N129: -13298.5 -> -13326.4
$ws.Cells.Item(137, 8).Value = 1315.0625  # H137: 1403.6786 -> 1315.0625
$ws.Cells.Item(137, 9).Value = 1375.8077  # I137: 1464.8695 -> 1375.8077
$ws.Cells.Item(137, 10).Value = 1051.8334  # J137: 1122.2 -> 1051.8334
$ws.Cells.Item(137, 11).Value = 4127.4231  # K137: 4394.6085 -> 4127.4231
$ws.Cells.Item(137, 12).Value = 3155.5002  # L137: 3366.6 -> 3155.5002
$ws.Cells.Item(137, 13).Value = -1577.4231  # M137: -1844.6085 -> -1577.4231
$ws.Cells.Item(137, 14).Value = -8255.5002  # N137: -8466.6 -> -8255.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(46, 8).Value = 3476  # H46: 3500 -> 3476
$ws.Cells.Item(46, 10).Value = 3952  # J46: 4000 -> 3952
$ws.Cells.Item(46, 12).Value = 3952  # L46: 4000 -> 3952
$ws.Cells.Item(46, 14).Value = -4590  # N46: -4638 -> -4590
$ws.Cells.Item(97, 8).Value = 29618.914  # H97: 28819.5 -> 29618.914
$ws.Cells.Item(97, 9).Value = 37746.703  # I97: 37748.184 -> 37746.703
$ws.Cells.Item(97, 10).Value = 2187.625  # J97: 2033.4445 -> 2187.625
$ws.Cells.Item(97, 11).Value = 37746.703  # K97: 37748.184 -> 37746.703
$ws.Cells.Item(97, 12).Value = 2187.625  # L97: 2033.4445 -> 2187.625
$ws.Cells.Item(97, 13).Value = -37250.703  # M97: -37252.184 -> -37250.703
$ws.Cells.Item(97, 14).Value = -3179.625  # N97: -3025.4445 -> -3179.625
$ws.Cells.Item(102, 8).Value = 102009.1  # H102: 127185.125 -> 102009.1
$ws.Cells.Item(102, 9).Value = 168698.33  # I102: 252395 -> 168698.33
$ws.Cells.Item(102, 11).Value = 168698.33  # K102: 252395 -> 168698.33
$ws.Cells.Item(102, 13).Value = -167076.33  # M102: -250773 -> -167076.33
$ws.Cells.Item(110, 8).Value = 34518080  # H110: 40044850 -> 34518080
$ws.Cells.Item(110, 9).Value = 41709160  # I110: 50055830 -> 41709160
$ws.Cells.Item(110, 10).Value = 886  # J110: 911 -> 886
$ws.Cells.Item(110, 11).Value = 41709160  # K110: 50055830 -> 41709160
$ws.Cells.Item(110, 12).Value = 886  # L110: 911 -> 886
$ws.Cells.Item(110, 13).Value = -41707115  # M110: -50053785 -> -41707115
$ws.Cells.Item(110, 14).Value = -4976  # N110: -5001 -> -4976
$ws.Cells.Item(132, 8).Value = 1379  # H132: 1506 -> 1379
$ws.Cells.Item(132, 9).Value = 1373.8937  # I132: 1506.3414 -> 1373.8937
$ws.Cells.Item(132, 11).Value = 4121.6811  # K132: 4519.0242 -> 4121.6811
$ws.Cells.Item(132, 13).Value = -1591.6811  # M132: -1989.0242 -> -1591.6811

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 41216.96  # H20: 50841.24 -> 41216.96
$ws.Cells.Item(20, 9).Value = 55880.684  # I20: 70570.47 -> 55880.684
$ws.Cells.Item(20, 10).Value = 1415.4286  # J20: 1518.1666 -> 1415.4286
$ws.Cells.Item(20, 11).Value = 55880.684  # K20: 70570.47 -> 55880.684
$ws.Cells.Item(20, 12).Value = 1415.4286  # L20: 1518.1666 -> 1415.4286
$ws.Cells.Item(20, 13).Value = -55633.684  # M20: -70323.47 -> -55633.684
$ws.Cells.Item(20, 14).Value = -1909.4286  # N20: -2012.1666 -> -1909.4286
$ws.Cells.Item(99, 8).Value = 1427.5714  # H99: 1345.2051 -> 1427.5714
$ws.Cells.Item(99, 9).Value = 1090.5555  # I99: 947.1539 -> 1090.5555
$ws.Cells.Item(99, 11).Value = 1090.5555  # K99: 947.1539 -> 1090.5555
$ws.Cells.Item(99, 13).Value = 407.4445000000001  # M99: 550.8461 -> 407.4445000000001
$ws.Cells.Item(107, 8).Value = 66698584  # H107: 66698588 -> 66698584
$ws.Cells.Item(134, 8).Value = 1888.2273  # H134: 1870.2683 -> 1888.2273
$ws.Cells.Item(134, 9).Value = 1586.0256  # I134: 1612.2894 -> 1586.0256
$ws.Cells.Item(134, 10).Value = 4245.4  # J134: 5138 -> 4245.4
$ws.Cells.Item(134, 11).Value = 4758.0768  # K134: 4836.8682 -> 4758.0768
$ws.Cells.Item(134, 12).Value = 12736.2  # L134: 15414 -> 12736.2
$ws.Cells.Item(134, 13).Value = -2223.0768  # M134: -2301.8682 -> -2223.0768
$ws.Cells.Item(134, 14).Value = -17806.2  # N134: -20484 -> -17806.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 27934.709  # H31: 30602.18 -> 27934.709
$ws.Cells.Item(31, 9).Value = 1665.64  # I31: 1744.8096 -> 1665.64
$ws.Cells.Item(31, 10).Value = 49825.6  # J31: 51498.9 -> 49825.6
$ws.Cells.Item(31, 11).Value = 1665.64  # K31: 1744.8096 -> 1665.64
$ws.Cells.Item(31, 12).Value = 49825.6  # L31: 51498.9 -> 49825.6
$ws.Cells.Item(31, 13).Value = -1370.64  # M31: -1449.8096 -> -1370.64
$ws.Cells.Item(31, 14).Value = -50415.6  # N31: -52088.9 -> -50415.6
$ws.Cells.Item(34, 8).Value = 27934.709  # H34: 30602.18 -> 27934.709
$ws.Cells.Item(34, 9).Value = 1665.64  # I34: 1744.8096 -> 1665.64
$ws.Cells.Item(34, 10).Value = 49825.6  # J34: 51498.9 -> 49825.6
$ws.Cells.Item(34, 11).Value = 1665.64  # K34: 1744.8096 -> 1665.64
$ws.Cells.Item(34, 12).Value = 49825.6  # L34: 51498.9 -> 49825.6
$ws.Cells.Item(34, 13).Value = -1463.64  # M34: -1542.8096 -> -1463.64
$ws.Cells.Item(34, 14).Value = -50229.6  # N34: -51902.9 -> -50229.6
$ws.Cells.Item(99, 8).Value = 3207.2354  # H99: 2989.5789 -> 3207.2354
$ws.Cells.Item(99, 9).Value = 4360  # I99: 3181.3333 -> 4360
$ws.Cells.Item(99, 10).Value = 2726.9167  # J99: 2817 -> 2726.9167
$ws.Cells.Item(99, 11).Value = 4360  # K99: 3181.3333 -> 4360
$ws.Cells.Item(99, 12).Value = 2726.9167  # L99: 2817 -> 2726.9167
$ws.Cells.Item(99, 13).Value = -2862  # M99: -1683.3333 -> -2862
$ws.Cells.Item(99, 14).Value = -5722.9167  # N99: -5813 -> -5722.9167
$ws.Cells.Item(122, 8).Value = 1398.4445  # H122: 1557.5333 -> 1398.4445
$ws.Cells.Item(122, 9).Value = 893.26666  # I122: 953.5 -> 893.26666
$ws.Cells.Item(122, 10).Value = 3924.3333  # J122: 10014 -> 3924.3333
$ws.Cells.Item(122, 11).Value = 2679.79998  # K122: 2860.5 -> 2679.79998
$ws.Cells.Item(122, 12).Value = 11772.9999  # L122: 30042 -> 11772.9999
$ws.Cells.Item(122, 13).Value = -229.7999799999998  # M122: -410.5 -> -229.7999799999998
$ws.Cells.Item(122, 14).Value = -16672.9999  # N122: -34942 -> -16672.9999
$ws.Cells.Item(126, 8).Value = 3207.2354  # H126: 2989.5789 -> 3207.2354
$ws.Cells.Item(126, 9).Value = 4360  # I126: 3181.3333 -> 4360
$ws.Cells.Item(126, 10).Value = 2726.9167  # J126: 2817 -> 2726.9167
$ws.Cells.Item(126, 11).Value = 13080  # K126: 9543.999899999999 -> 13080
$ws.Cells.Item(126, 12).Value = 8180.750100000001  # L126: 8451 -> 8180.750100000001
$ws.Cells.Item(126, 13).Value = -10610  # M126: -7073.999899999999 -> -10610
$ws.Cells.Item(126, 14).Value = -13120.7501  # N126: -13391 -> -13120.7501
$ws.Cells.Item(134, 8).Value = 1512.25  # H134: 1351.3636 -> 1512.25
$ws.Cells.Item(134, 9).Value = 1472.6923  # I134: 1238.5625 -> 1472.6923
$ws.Cells.Item(134, 10).Value = 1585.7142  # J134: 1652.1666 -> 1585.7142
$ws.Cells.Item(134, 11).Value = 4418.0769  # K134: 3715.6875 -> 4418.0769
$ws.Cells.Item(134, 12).Value = 4757.142599999999  # L134: 4956.4998 -> 4757.142599999999
$ws.Cells.Item(134, 13).Value = -1883.0769  # M134: -1180.6875 -> -1883.0769
$ws.Cells.Item(134, 14).Value = -9827.142599999999  # N134: -10026.4998 -> -9827.142599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 2000  # H3: 1103.3334 -> 2000
$ws.Cells.Item(3, 9).Value = 2000  # I3: 1103.3334 -> 2000
$ws.Cells.Item(3, 11).Value = 6000  # K3: 3310.0002 -> 6000
$ws.Cells.Item(3, 13).Value = -5888  # M3: -3198.0002 -> -5888
$ws.Cells.Item(131, 8).Value = 764.96  # H131: 758.52 -> 764.96
$ws.Cells.Item(131, 10).Value = 796.3511  # J131: 789.5 -> 796.3511
$ws.Cells.Item(131, 12).Value = 2389.0533  # L131: 2368.5 -> 2389.0533
$ws.Cells.Item(131, 14).Value = -12469.0533  # N131: -12448.5 -> -12469.0533

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 1835555.9  # H126: 1541261.1 -> 1835555.9
$ws.Cells.Item(126, 9).Value = 3145039.8  # I126: 2503125.2 -> 3145039.8
$ws.Cells.Item(126, 11).Value = 9435119.399999999  # K126: 7509375.600000001 -> 9435119.399999999
$ws.Cells.Item(126, 13).Value = -9432649.399999999  # M126: -7506905.600000001 -> -9432649.399999999
$ws.Cells.Item(132, 8).Value = 2156.5283  # H132: 2374.7112 -> 2156.5283
$ws.Cells.Item(132, 9).Value = 2082.0852  # I132: 2273.756 -> 2082.0852
$ws.Cells.Item(132, 10).Value = 2739.6667  # J132: 3409.5 -> 2739.6667
$ws.Cells.Item(132, 11).Value = 6246.2556  # K132: 6821.268 -> 6246.2556
$ws.Cells.Item(132, 12).Value = 8219.000100000001  # L132: 10228.5 -> 8219.000100000001
$ws.Cells.Item(132, 13).Value = -3716.2556  # M132: -4291.268 -> -3716.2556
$ws.Cells.Item(132, 14).Value = -13279.0001  # N132: -15288.5 -> -13279.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 93072.73  # H40: 251875 -> 93072.73
$ws.Cells.Item(40, 9).Value = 334600  # I40: 1000000 -> 334600
$ws.Cells.Item(40, 11).Value = 334600  # K40: 1000000 -> 334600
$ws.Cells.Item(40, 13).Value = -334464  # M40: -999864 -> -334464
$ws.Cells.Item(46, 8).Value = 1013035.6  # H46: 1125519.8 -> 1013035.6
$ws.Cells.Item(46, 10).Value = 1125540.6  # J46: 1266148.5 -> 1125540.6
$ws.Cells.Item(46, 12).Value = 1125540.6  # L46: 1266148.5 -> 1125540.6
$ws.Cells.Item(46, 14).Value = -1125916.6  # N46: -1266524.5 -> -1125916.6
$ws.Cells.Item(61, 8).Value = 1584.9565  # H61: 1695 -> 1584.9565
$ws.Cells.Item(61, 9).Value = 1613.75  # I61: 1781.1 -> 1613.75
$ws.Cells.Item(61, 10).Value = 1553.5454  # J61: 1608.9 -> 1553.5454
$ws.Cells.Item(61, 11).Value = 1613.75  # K61: 1781.1 -> 1613.75
$ws.Cells.Item(61, 12).Value = 1553.5454  # L61: 1608.9 -> 1553.5454
$ws.Cells.Item(61, 13).Value = -1411.75  # M61: -1579.1 -> -1411.75
$ws.Cells.Item(61, 14).Value = -1957.5454  # N61: -2012.9 -> -1957.5454
$ws.Cells.Item(93, 8).Value = 2516.3157  # H93: 1815.5186 -> 2516.3157
$ws.Cells.Item(93, 9).Value = 2661.3333  # I93: 1768.7778 -> 2661.3333
$ws.Cells.Item(93, 10).Value = 2267.7144  # J93: 1909 -> 2267.7144
$ws.Cells.Item(93, 11).Value = 2661.3333  # K93: 1768.7778 -> 2661.3333
$ws.Cells.Item(93, 12).Value = 2267.7144  # L93: 1909 -> 2267.7144
$ws.Cells.Item(93, 13).Value = -1413.3333  # M93: -520.7778000000001 -> -1413.3333
$ws.Cells.Item(93, 14).Value = -4763.7144  # N93: -4405 -> -4763.7144
$ws.Cells.Item(113, 8).Value = 1584.9565  # H113: 1695 -> 1584.9565
$ws.Cells.Item(113, 9).Value = 1613.75  # I113: 1781.1 -> 1613.75
$ws.Cells.Item(113, 10).Value = 1553.5454  # J113: 1608.9 -> 1553.5454
$ws.Cells.Item(113, 11).Value = 1613.75  # K113: 1781.1 -> 1613.75
$ws.Cells.Item(113, 12).Value = 1553.5454  # L113: 1608.9 -> 1553.5454
$ws.Cells.Item(113, 13).Value = 556.25  # M113: 388.9000000000001 -> 556.25
$ws.Cells.Item(113, 14).Value = -5893.5454  # N113: -5948.9 -> -5893.5454
$ws.Cells.Item(132, 8).Value = 4408.5713  # H132: 4622.222 -> 4408.5713
$ws.Cells.Item(132, 9).Value = 4354  # I132: 4622.222 -> 4354
$ws.Cells.Item(132, 10).Value = 5500  # J132: 0 -> 5500
$ws.Cells.Item(132, 11).Value = 13062  # K132: 13866.666 -> 13062
$ws.Cells.Item(132, 12).Value = 16500  # L132: 0 -> 16500
$ws.Cells.Item(132, 13).Value = -10532  # M132: -11336.666 -> -10532
$ws.Cells.Item(132, 14).Value = -21560  # N132: (blank) -> -21560
$ws.Cells.Item(136, 8).Value = 2505.611  # H136: 2517.7058 -> 2505.611
$ws.Cells.Item(136, 9).Value = 2321.5  # I136: 2335.7856 -> 2321.5
$ws.Cells.Item(136, 10).Value = 3150  # J136: 3366.6667 -> 3150
$ws.Cells.Item(136, 11).Value = 6964.5  # K136: 7007.3568 -> 6964.5
$ws.Cells.Item(136, 12).Value = 9450  # L136: 10100.0001 -> 9450
$ws.Cells.Item(136, 13).Value = -4414.5  # M136: -4457.3568 -> -4414.5
$ws.Cells.Item(136, 14).Value = -14550  # N136: -15200.0001 -> -14550
$ws.Cells.Item(141, 8).Value = 0  # H141: 66043.336 -> 0
$ws.Cells.Item(141, 10).Value = 0  # J141: 66043.336 -> 0
$ws.Cells.Item(141, 12).Value = 0  # L141: 66043.336 -> 0
$ws.Cells.Item(141, 14).ClearContents()  # N141: remove (was -76403.336)

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(32, 8).Value = 9896.25  # H32: 11762.917 -> 9896.25
$ws.Cells.Item(32, 9).Value = 2651.6667  # I32: 5555 -> 2651.6667
$ws.Cells.Item(32, 10).Value = 12311.111  # J32: 12327.272 -> 12311.111
$ws.Cells.Item(32, 11).Value = 2651.6667  # K32: 5555 -> 2651.6667
$ws.Cells.Item(32, 12).Value = 12311.111  # L32: 12327.272 -> 12311.111
$ws.Cells.Item(32, 13).Value = -2334.6667  # M32: -5238 -> -2334.6667
$ws.Cells.Item(32, 14).Value = -12945.111  # N32: -12961.272 -> -12945.111
$ws.Cells.Item(122, 8).Value = 1558.381  # H122: 1708.1666 -> 1558.381
$ws.Cells.Item(122, 9).Value = 1263.0667  # I122: 1413.9166 -> 1263.0667
$ws.Cells.Item(122, 11).Value = 3789.2001  # K122: 4241.7498 -> 3789.2001
$ws.Cells.Item(122, 13).Value = -1339.2001  # M122: -1791.7498 -> -1339.2001
$ws.Cells.Item(126, 8).Value = 1408.5  # H126: 1534.9166 -> 1408.5
$ws.Cells.Item(126, 9).Value = 1235  # I126: 1352 -> 1235
$ws.Cells.Item(126, 11).Value = 3705  # K126: 4056 -> 3705
$ws.Cells.Item(126, 13).Value = -1235  # M126: -1586 -> -1235
